# The log sheet already contains two "runs" of data in rows 7:8 and 9:10
# (label "a1" / "b2" in column A, followed by the same set of simulation
# results in columns B:CV). This change appends one more run by
# duplicating the most recent run (rows 9:10) into new rows 11:12,
# extending the used range from A1:CV10 to A1:CV12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A9:CV10")
$destination = $ws.Range("A11:CV12")

$source.Copy($destination)
